$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph: "LISTEN ON (for APPLE, with triggered moves) and "
# Becomes 5 runs:
#   "LISTEN " / "IOS" / " (for APPLE, with triggered moves" / ") " /
#   "6 sec moving cycle (3sec twisted + 3sec at home)"
# ------------------------------------------------------------------

# Shrink the existing run down to just "LISTEN " first.
$d.Content.Find.Execute(
    "LISTEN ON (for APPLE, with triggered moves) and ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "LISTEN ",
    2)

# Turning TrackRevisions on while inserting the remaining text keeps each
# insertion as its own run (rather than silently coalescing into the
# neighbouring run of identical formatting); AcceptAllRevisions() then
# folds the insertions into the body without leaving any revision markup.
$d.TrackRevisions = $true

$listenPara = $d.Paragraphs.Item(39)
$r = $listenPara.Range
$insPoint = $r.End - 1

$ins1 = $d.Range($insPoint, $insPoint)
$ins1.InsertAfter("IOS")
$ins1.Font.Name = "Century Gothic"

$ins2 = $d.Range($ins1.End, $ins1.End)
$ins2.InsertAfter(" (for APPLE, with triggered moves")
$ins2.Font.Name = "Century Gothic"

$ins3 = $d.Range($ins2.End, $ins2.End)
$ins3.InsertAfter(") ")
$ins3.Font.Name = "Century Gothic"

$ins4 = $d.Range($ins3.End, $ins3.End)
$ins4.InsertAfter("6 sec moving cycle (3sec twisted + 3sec at home)")
$ins4.Font.Name = "Century Gothic"

# ------------------------------------------------------------------
# Paragraph: "LISTEN FULL for verbose mode, no moves triggered"
# Gains a second run: ". Update every .3 seconds"
# ------------------------------------------------------------------
$fullPara = $d.Paragraphs.Item(40)
$r2 = $fullPara.Range
$insPoint2 = $r2.End - 1

$ins5 = $d.Range($insPoint2, $insPoint2)
$ins5.InsertAfter(". Update every .3 seconds")
$ins5.Font.Name = "Century Gothic"

$d.TrackRevisions = $false
$d.AcceptAllRevisions()

# ------------------------------------------------------------------
# Two new trailing paragraphs, matching the formatting of the
# paragraph they are appended after.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last
$p1.Range.Text = "Negative values are on the right, positives on the left."

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Run a MICBAL in a quiet room to give the ears a baseline. Value stored in flash."
